$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "2024-03-28") {
                $sh.TextFrame.TextRange.Text = "2024-03-29"
            }
        }
    }
}

# Update the slide master's date placeholder
$m = $p.SlideMaster
Update-DateShape $m.Shapes

# Update every slide layout's date placeholder
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    Update-DateShape $layout.Shapes
}
